$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 518.6316
$ws.Range("I18").Value = 347.125
$ws.Range("J18").Value = 1433.3334
$ws.Range("K18").Value = 347.125
$ws.Range("L18").Value = 1433.3334
$ws.Range("M18").Value = -63.125
$ws.Range("N18").Value = -2001.3334

$ws.Range("H43").Value = 625.25
$ws.Range("I43").Value = 583.3333
$ws.Range("J43").Value = 751
$ws.Range("K43").Value = 583.3333
$ws.Range("L43").Value = 751
$ws.Range("M43").Value = -514.3333
$ws.Range("N43").Value = -889

$ws.Range("H129").Value = 999.59375
$ws.Range("I129").Value = 321.75
$ws.Range("J129").Value = 1096.4286
$ws.Range("K129").Value = 965.25
$ws.Range("L129").Value = 3289.2858
$ws.Range("M129").Value = 4034.75
$ws.Range("N129").Value = -13289.2858

$ws.Range("H132").Value = 1702.8182
$ws.Range("I132").Value = 1409.4517
$ws.Range("J132").Value = 6250
$ws.Range("K132").Value = 4228.355100000001
$ws.Range("L132").Value = 18750
$ws.Range("M132").Value = -1698.355100000001
$ws.Range("N132").Value = -23810

$ws.Range("H137").Value = 2589.111
$ws.Range("I137").Value = 2551.7576
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 7655.2728
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -5105.2728
$ws.Range("N137").Value = -14100


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29922.838
$ws.Range("I32").Value = 32978.21
$ws.Range("J32").Value = 4716
$ws.Range("K32").Value = 32978.21
$ws.Range("L32").Value = 4716
$ws.Range("M32").Value = -32691.21
$ws.Range("N32").Value = -5290

$ws.Range("H61").Value = 7322.3423
$ws.Range("I61").Value = 5313.44
$ws.Range("J61").Value = 11185.615
$ws.Range("K61").Value = 5313.44
$ws.Range("L61").Value = 11185.615
$ws.Range("M61").Value = -5101.44
$ws.Range("N61").Value = -11609.615

$ws.Range("H63").Value = 3427.1428
$ws.Range("I63").Value = 2670
$ws.Range("J63").Value = 3995
$ws.Range("K63").Value = 2670
$ws.Range("L63").Value = 3995
$ws.Range("M63").Value = -1984
$ws.Range("N63").Value = -5367

$ws.Range("H66").Value = 3427.1428
$ws.Range("I66").Value = 2670
$ws.Range("J66").Value = 3995
$ws.Range("K66").Value = 13350
$ws.Range("L66").Value = 19975
$ws.Range("M66").Value = -9918
$ws.Range("N66").Value = -26839

$ws.Range("H74").Value = 2262.4062
$ws.Range("I74").Value = 2042.76
$ws.Range("J74").Value = 3046.8572
$ws.Range("K74").Value = 2042.76
$ws.Range("L74").Value = 3046.8572
$ws.Range("M74").Value = -1168.76
$ws.Range("N74").Value = -4794.8572

$ws.Range("H77").Value = 2262.4062
$ws.Range("I77").Value = 2042.76
$ws.Range("J77").Value = 3046.8572
$ws.Range("K77").Value = 10213.8
$ws.Range("L77").Value = 15234.286
$ws.Range("M77").Value = -5845.799999999999
$ws.Range("N77").Value = -23970.286

$ws.Range("H122").Value = 1926
$ws.Range("I122").Value = 1977.875
$ws.Range("J122").Value = 1760
$ws.Range("K122").Value = 5933.625
$ws.Range("L122").Value = 5280
$ws.Range("M122").Value = -3483.625
$ws.Range("N122").Value = -10180

$ws.Range("H136").Value = 7322.3423
$ws.Range("I136").Value = 5313.44
$ws.Range("J136").Value = 11185.615
$ws.Range("K136").Value = 15940.32
$ws.Range("L136").Value = 33556.845
$ws.Range("M136").Value = -13390.32
$ws.Range("N136").Value = -38656.845


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 70000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70000
$ws.Range("N35").Value = -70620

$ws.Range("H80").Value = 113.53333
$ws.Range("I80").Value = 87
$ws.Range("J80").Value = 120.166664
$ws.Range("K80").Value = 87
$ws.Range("L80").Value = 120.166664
$ws.Range("M80").Value = 911
$ws.Range("N80").Value = -2116.166664

$ws.Range("H83").Value = 113.53333
$ws.Range("I83").Value = 87
$ws.Range("J83").Value = 120.166664
$ws.Range("K83").Value = 435
$ws.Range("L83").Value = 600.83332
$ws.Range("M83").Value = 4557
$ws.Range("N83").Value = -10584.83332

$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4996


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7414.55
$ws.Range("I31").Value = 7176.4243
$ws.Range("J31").Value = 8537.143
$ws.Range("K31").Value = 7176.4243
$ws.Range("L31").Value = 8537.143
$ws.Range("M31").Value = -6881.4243
$ws.Range("N31").Value = -9127.143

$ws.Range("H34").Value = 7414.55
$ws.Range("I34").Value = 7176.4243
$ws.Range("J34").Value = 8537.143
$ws.Range("K34").Value = 7176.4243
$ws.Range("L34").Value = 8537.143
$ws.Range("M34").Value = -6974.4243
$ws.Range("N34").Value = -8941.143

$ws.Range("H132").Value = 5599.3057
$ws.Range("I132").Value = 6379.32
$ws.Range("J132").Value = 3826.5454
$ws.Range("K132").Value = 19137.96
$ws.Range("L132").Value = 11479.6362
$ws.Range("M132").Value = -16607.96
$ws.Range("N132").Value = -16539.6362

$ws.Range("H134").Value = 2547.795
$ws.Range("I134").Value = 2067.36
$ws.Range("J134").Value = 3405.7144
$ws.Range("K134").Value = 6202.08
$ws.Range("L134").Value = 10217.1432
$ws.Range("M134").Value = -3667.08
$ws.Range("N134").Value = -15287.1432


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.142857
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 34.962963
$ws.Range("K2").Value = 72
$ws.Range("L2").Value = 209.777778
$ws.Range("M2").Value = 41
$ws.Range("N2").Value = -435.777778

$ws.Range("H114").Value = 760.55554
$ws.Range("I114").Value = 428.5
$ws.Range("J114").Value = 1026.2
$ws.Range("K114").Value = 1285.5
$ws.Range("L114").Value = 3078.6
$ws.Range("M114").Value = 1968.5
$ws.Range("N114").Value = -9586.6

$ws.Range("H119").Value = 6321.5
$ws.Range("I119").Value = 464.5
$ws.Range("J119").Value = 9250
$ws.Range("K119").Value = 1393.5
$ws.Range("L119").Value = 27750
$ws.Range("M119").Value = 3444.5
$ws.Range("N119").Value = -37426

$ws.Range("H141").Value = 3628.3333
$ws.Range("I141").Value = 1948.5714
$ws.Range("J141").Value = 5980
$ws.Range("K141").Value = 5845.7142
$ws.Range("L141").Value = 17940
$ws.Range("M141").Value = -665.7142000000003
$ws.Range("N141").Value = -28300


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 15587.5
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 15587.5
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15587.5
$ws.Range("N52").Value = -16105.5

$ws.Range("H132").Value = 7591.9473
$ws.Range("I132").Value = 11291.6
$ws.Range("J132").Value = 3481.2222
$ws.Range("K132").Value = 33874.8
$ws.Range("L132").Value = 10443.6666
$ws.Range("M132").Value = -31344.8
$ws.Range("N132").Value = -15503.6666


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 23196
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 23196
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 23196
$ws.Range("N25").Value = -23656

$ws.Range("H132").Value = 5919.357
$ws.Range("I132").Value = 6200.6665
$ws.Range("J132").Value = 5708.375
$ws.Range("K132").Value = 18601.9995
$ws.Range("L132").Value = 17125.125
$ws.Range("M132").Value = -16071.9995
$ws.Range("N132").Value = -22185.125


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4000
$ws.Range("N17").Value = -4344

$ws.Range("H81").Value = 9093869
$ws.Range("I81").Value = 1620
$ws.Range("J81").Value = 15388503
$ws.Range("K81").Value = 3240
$ws.Range("L81").Value = 30777006
$ws.Range("M81").Value = -2179
$ws.Range("N81").Value = -30779128

$ws.Range("H84").Value = 9093869
$ws.Range("I84").Value = 1620
$ws.Range("J84").Value = 15388503
$ws.Range("K84").Value = 16200
$ws.Range("L84").Value = 153885030
$ws.Range("M84").Value = -10896
$ws.Range("N84").Value = -153895638


# Remove cells that are deleted entirely in the target state
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M25").Value = $null
